$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the top of the data (row 2), pushing every
# existing data row (and its hyperlink) down by one. This alone turns the
# old row 196 (07-08-2025) into the new row 197, matching the target.
$ws.Rows("2:2").Insert()

# Populate the newly-inserted row 2 with the latest price entry.
$ws.Cells.Item(2, 1).Value2 = "18-02-2026"
$ws.Cells.Item(2, 2).Value2 = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value2 = "IE07"
$ws.Cells.Item(2, 4).Value2 = 320.45

# Column E holds "01-02-2026" which Excel's smart-entry would otherwise
# read as an ambiguous m/d date; force text so it stays a literal string
# like the rest of the sheet, then restore the General format so the
# cell's appearance still matches its neighbours.
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value2 = "01-02-2026"
$ws.Cells.Item(2, 5).NumberFormat = "General"

$ws.Cells.Item(2, 6).Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-02-2026.pdf"

# Re-create the hyperlink on F2 (Insert() only carried over the old F2's
# link target; repoint it at the new circular's PDF).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-02-2026.pdf")
